$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (row2 / row3) had their trial-parameter values (cols B:J)
# swapped, while column A (trial index) stayed as-is.
$row2Values = $ws.Range("B2:J2").Value2
$row3Values = $ws.Range("B3:J3").Value2

$ws.Range("B2:J2").Value2 = $row3Values
$ws.Range("B3:J3").Value2 = $row2Values
